$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E for "ram" - copy the header formatting from D1 (bold/border/center)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "ram"

$ws.Range("E2").Value = 12
$ws.Range("E3").Value = "M"
$ws.Range("E4").Value = 400
$ws.Range("E5").Value = 200
$ws.Range("E6").Value = 1.63
$ws.Range("E7").Value = 1.67
$ws.Range("E8").Value = 16.93703187925778
$ws.Range("E9").Value = 16.85252249991036
